$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last five comments (old rows 21-25) were dropped entirely, which
# shifts every row below the header up so the sheet now ends at row 20.
$ws.Rows("21:25").Delete()

# Column B holds literal "DD-MM-YYYY" text (not real dates), so force
# a text format first or Excel would silently coerce it to a date serial.
$ws.Range("B2:B20").NumberFormat = "@"

# Rewrite rows 2-20 with the realigned/updated data (text shifted up one
# row relative to author/date/concept, plus several value corrections).
# Row 2
$ws.Range("A2").Value = "Son movidas, gana massa baja el dólar para que su campaña sea excelente"
$ws.Range("B2").Value = "10-10-2023"
$ws.Range("C2").Value = "franquiito.02"
$ws.Range("D2").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E2").Value = "alegría"
$ws.Range("F2").Value = "N/E"
$ws.Range("G2").Value = "N/E"
$ws.Range("H2").Value = "N/E"

# Row 3
$ws.Range("A3").Value = "100% dee acuerdo"
$ws.Range("B3").Value = "09-10-2023"
$ws.Range("C3").Value = "fernando_taboas"
$ws.Range("D3").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E3").Value = "alegría"
$ws.Range("F3").Value = "N/E"
$ws.Range("G3").Value = "N/E"
$ws.Range("H3").Value = "N/E"

# Row 4
$ws.Range("A4").Value = "Es totalmente lógico lo que dicen...pero a cuánto estaría si estuviesen los liberales o jxc??? A menos, seguro??? O lo hubiesen liberado y hoy estría a 3000?"
$ws.Range("B4").Value = "10-10-2023"
$ws.Range("C4").Value = "claudiogabrielvillafanie"
$ws.Range("D4").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E4").Value = "tristeza"
$ws.Range("F4").Value = "N/E"
$ws.Range("G4").Value = "N/E"
$ws.Range("H4").Value = "N/E"

# Row 5
$ws.Range("A5").Value = "Lo mal que está Argentina en lo económico y social, es la herencia dejada por el gobierno corrupto de Mauricio Macri, difícil fue para el gobierno de Alberto Fernández, revertir en tan poco tiempo el desastre planificado por la derecha política y Empresarial Argentina.`nEs de esperar que el Gobierno del Presidente Massa, pueda tomar las medidas económicas necesarias,  para corregir esta pesada situación que afecta al pueblo Argentino."
$ws.Range("B5").Value = "10-10-2023"
$ws.Range("C5").Value = "jg.cuevasc020653"
$ws.Range("D5").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E5").Value = "enojo"
$ws.Range("F5").Value = "N/E"
$ws.Range("G5").Value = "N/E"
$ws.Range("H5").Value = "N/E"

# Row 6
$ws.Range("A6").Value = "ah pero macri"
$ws.Range("B6").Value = "10-10-2023"
$ws.Range("C6").Value = "_scared_5"
$ws.Range("D6").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E6").Value = "alegría"
$ws.Range("F6").Value = "N/E"
$ws.Range("G6").Value = "N/E"
$ws.Range("H6").Value = "N/E"

# Row 7
$ws.Range("A7").Value = "1050, antes de que gane Alberto yo me agarraba de los pelos, era obvio lo que iba a pasar pero los K tiraron de la soga y ya se está por romperrr"
$ws.Range("B7").Value = "10-10-2023"
$ws.Range("C7").Value = "martin_zingoni"
$ws.Range("D7").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E7").Value = "sorpresa"
$ws.Range("F7").Value = "N/E"
$ws.Range("G7").Value = "N/E"
$ws.Range("H7").Value = "N/E"

# Row 8
$ws.Range("A8").Value = "Vergüenza es votar a milei...."
$ws.Range("B8").Value = "10-10-2023"
$ws.Range("C8").Value = "hora.vizueta"
$ws.Range("D8").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E8").Value = "asco"
$ws.Range("F8").Value = "N/E"
$ws.Range("G8").Value = "N/E"
$ws.Range("H8").Value = "N/E"

# Row 9
$ws.Range("A9").Value = "#MassaPresidente"
$ws.Range("B9").Value = "09-10-2023"
$ws.Range("C9").Value = "drgustavovaldez"
$ws.Range("D9").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E9").Value = "alegría"
$ws.Range("F9").Value = "N/E"
$ws.Range("G9").Value = "N/E"
$ws.Range("H9").Value = "N/E"

# Row 10
$ws.Range("A10").Value = "Un desastre"
$ws.Range("B10").Value = "10-10-2023"
$ws.Range("C10").Value = "silvina.harris"
$ws.Range("D10").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E10").Value = "alegría"
$ws.Range("F10").Value = "N/E"
$ws.Range("G10").Value = "N/E"
$ws.Range("H10").Value = "N/E"

# Row 11
$ws.Range("A11").Value = "Jajaja Jajaja"
$ws.Range("B11").Value = "10-10-2023"
$ws.Range("C11").Value = "daro.gonzalez"
$ws.Range("D11").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E11").Value = "alegría"
$ws.Range("F11").Value = "N/E"
$ws.Range("G11").Value = "N/E"
$ws.Range("H11").Value = "N/E"

# Row 12
$ws.Range("A12").Value = "Ustedes dan vergüenza"
$ws.Range("B12").Value = "10-10-2023"
$ws.Range("C12").Value = "zarpadomal1959"
$ws.Range("D12").Value = "https://www.threads.net/@liberalesargentinaok/post/CyMS-ONuuaU"
$ws.Range("E12").Value = "asco"
$ws.Range("F12").Value = "N/E"
$ws.Range("G12").Value = "N/E"
$ws.Range("H12").Value = "N/E"

# Row 13
$ws.Range("A13").Value = "Si quieren terminar con el Kirchnerismo, Juntos por el Cambio debe apoyar a Javier Milei y el 19 de Noviembre terminan arrasando. Es simple. No jodan."
$ws.Range("B13").Value = "23-10-2023"
$ws.Range("C13").Value = "lanatappt"
$ws.Range("D13").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E13").Value = "enojo"
$ws.Range("F13").Value = "N/E"
$ws.Range("G13").Value = "N/E"
$ws.Range("H13").Value = "N/E"

# Row 14
$ws.Range("A14").Value = "La anti ética dixit...🤦‍♀️"
$ws.Range("B14").Value = "26-10-2023"
$ws.Range("C14").Value = "veronicacepedaph"
$ws.Range("D14").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E14").Value = "sorpresa"
$ws.Range("F14").Value = "N/E"
$ws.Range("G14").Value = "N/E"
$ws.Range("H14").Value = "N/E"

# Row 15
$ws.Range("A15").Value = "Te acordás cuándo estabas al frente de revista veintitrés? Te cabía Chávez, fidel y le dabas chirlos a Macri mientras fumabas en la ducha?"
$ws.Range("B15").Value = "24-10-2023"
$ws.Range("C15").Value = "leover.ok"
$ws.Range("D15").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E15").Value = "alegría"
$ws.Range("F15").Value = "N/E"
$ws.Range("G15").Value = "N/E"
$ws.Range("H15").Value = "N/E"

# Row 16
$ws.Range("A16").Value = "No hay otra posibilidad sino Ezeiza"
$ws.Range("B16").Value = "25-10-2023"
$ws.Range("C16").Value = "mariasusana_sartori"
$ws.Range("D16").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E16").Value = "sorpresa"
$ws.Range("F16").Value = "N/E"
$ws.Range("G16").Value = "N/E"
$ws.Range("H16").Value = "N/E"

# Row 17
$ws.Range("A17").Value = "Que Antipatria y Facho que sos Lanata, remítete a decir tu opinión ,No a decirle a la gente lo que tiene que hacer."
$ws.Range("B17").Value = "24-10-2023"
$ws.Range("C17").Value = "kristellrosario"
$ws.Range("D17").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E17").Value = "enojo"
$ws.Range("F17").Value = "N/E"
$ws.Range("G17").Value = "N/E"
$ws.Range("H17").Value = "N/E"

# Row 18
$ws.Range("A18").Value = "Mediquese"
$ws.Range("B18").Value = "24-10-2023"
$ws.Range("C18").Value = "_g_a_b_r_i_e_a"
$ws.Range("D18").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E18").Value = "enojo"
$ws.Range("F18").Value = "N/E"
$ws.Range("G18").Value = "N/E"
$ws.Range("H18").Value = "N/E"

# Row 19
$ws.Range("A19").Value = "Coherencia y dignidad nunca no?"
$ws.Range("B19").Value = "25-10-2023"
$ws.Range("C19").Value = "ivanburaok"
$ws.Range("D19").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E19").Value = "enojo"
$ws.Range("F19").Value = "N/E"
$ws.Range("G19").Value = "N/E"
$ws.Range("H19").Value = "N/E"

# Row 20
$ws.Range("A20").Value = "NOAMILEI"
$ws.Range("B20").Value = "24-10-2023"
$ws.Range("C20").Value = "florviterbo"
$ws.Range("D20").Value = "https://www.threads.net/@lanatappt/post/CyvjalRg-vx"
$ws.Range("E20").Value = "alegría"
$ws.Range("F20").Value = "N/E"
$ws.Range("G20").Value = "N/E"
$ws.Range("H20").Value = "N/E"
